$p = $ppt.ActivePresentation

# --- Slide 3: merge the "the " run with the following "numerical indicator..." run ---
# (Fixes an accidental mid-sentence run split; visible text is unchanged.)
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(5)
$tr3 = $shape3.TextFrame.TextRange
$sub3 = $tr3.Characters(30, 124)
$sub3.Text = "the numerical indicator of the learning rate (based on calculating the similarity of the proposed keyboard with the default)"

# --- Slide 5: nudge the two result screenshots slightly to the left ---
$s5 = $p.Slides.Item(5)
$pic1 = $s5.Shapes.Item(5)
$pic1.Left = 52.81259842519685
$pic2 = $s5.Shapes.Item(6)
$pic2.Left = 971.339842519685

# --- Slide 6: fix the typo "which generate algorithm" -> "which was generated by algorithm" ---
$s6 = $p.Slides.Item(6)
$shape6 = $s6.Shapes.Item(2)
$tr6 = $shape6.TextFrame.TextRange
$tr6.Text = "Keyboard_prefix_10 is the best layout "
[void]$tr6.InsertAfter("which was generated by ")
[void]$tr6.InsertAfter("algorithm")
